$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (like "1.015" or "27.856.03")
# are written back as TEXT, matching the original inlineStr cell type,
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.856.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.70"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4717"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3918"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08058"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.021"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.03"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.918.28"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.965"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.099"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.016"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06751"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.886.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.505"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.135.62"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.48"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.09"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.094"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.518"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.96"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9739"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09474"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.438"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.646"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.349"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06138"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02263"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.221"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5986"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.008"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1892"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.28"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.265"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5691"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.14"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.411"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.934"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06922"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.67"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.072"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.01%  "
